$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 814.0909
$ws.Cells.Item(6, 9).Value = 264.72726
$ws.Cells.Item(6, 10).Value = 1363.4546
$ws.Cells.Item(6, 11).Value = 794.18178
$ws.Cells.Item(6, 12).Value = 4090.3638
$ws.Cells.Item(6, 13).Value = -682.18178
$ws.Cells.Item(6, 14).Value = -4314.3638
$ws.Cells.Item(17, 8).Value = 3342530.2
$ws.Cells.Item(17, 10).Value = 3413644.5
$ws.Cells.Item(17, 12).Value = 10240933.5
$ws.Cells.Item(17, 14).Value = -10241269.5
$ws.Cells.Item(43, 8).Value = 1465.4546
$ws.Cells.Item(43, 9).Value = 993.6667
$ws.Cells.Item(43, 10).Value = 1642.375
$ws.Cells.Item(43, 11).Value = 993.6667
$ws.Cells.Item(43, 12).Value = 1642.375
$ws.Cells.Item(43, 13).Value = -924.6667
$ws.Cells.Item(43, 14).Value = -1780.375
$ws.Cells.Item(100, 8).Value = 16668085
$ws.Cells.Item(100, 9).Value = 27779058
$ws.Cells.Item(100, 10).Value = 1624
$ws.Cells.Item(100, 11).Value = 27779058
$ws.Cells.Item(100, 12).Value = 1624
$ws.Cells.Item(100, 13).Value = -27778517
$ws.Cells.Item(100, 14).Value = -2706
$ws.Cells.Item(129, 8).Value = 859.4842
$ws.Cells.Item(129, 9).Value = 437.42105
$ws.Cells.Item(129, 10).Value = 965
$ws.Cells.Item(129, 11).Value = 1312.26315
$ws.Cells.Item(129, 12).Value = 2895
$ws.Cells.Item(129, 13).Value = 3687.73685
$ws.Cells.Item(129, 14).Value = -12895
$ws.Cells.Item(132, 8).Value = 827.0323
$ws.Cells.Item(132, 9).Value = 820.7143
$ws.Cells.Item(132, 11).Value = 2462.1429
$ws.Cells.Item(132, 13).Value = 67.85710000000017
$ws.Cells.Item(137, 8).Value = 1232.8286
$ws.Cells.Item(137, 9).Value = 1232
$ws.Cells.Item(137, 10).Value = 1237.8
$ws.Cells.Item(137, 11).Value = 3696
$ws.Cells.Item(137, 12).Value = 3713.4
$ws.Cells.Item(137, 13).Value = -1146
$ws.Cells.Item(137, 14).Value = -8813.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 14331.333
$ws.Cells.Item(13, 10).Value = 14331.333
$ws.Cells.Item(13, 12).Value = 14331.333
$ws.Cells.Item(13, 14).Value = -14619.333
$ws.Cells.Item(32, 8).Value = 4230.864
$ws.Cells.Item(32, 9).Value = 3014.1865
$ws.Cells.Item(32, 10).Value = 14485.714
$ws.Cells.Item(32, 11).Value = 3014.1865
$ws.Cells.Item(32, 12).Value = 14485.714
$ws.Cells.Item(32, 13).Value = -2727.1865
$ws.Cells.Item(32, 14).Value = -15059.714
$ws.Cells.Item(45, 8).Value = 15303.571
$ws.Cells.Item(45, 9).Value = 26103
$ws.Cells.Item(45, 10).Value = 904.3333
$ws.Cells.Item(45, 11).Value = 26103
$ws.Cells.Item(45, 12).Value = 904.3333
$ws.Cells.Item(45, 13).Value = -25726
$ws.Cells.Item(45, 14).Value = -1658.3333
$ws.Cells.Item(61, 8).Value = 9537.643
$ws.Cells.Item(61, 9).Value = 11554.454
$ws.Cells.Item(61, 11).Value = 11554.454
$ws.Cells.Item(61, 13).Value = -11342.454
$ws.Cells.Item(74, 8).Value = 1500.1578
$ws.Cells.Item(74, 9).Value = 1444.6111
$ws.Cells.Item(74, 10).Value = 2500
$ws.Cells.Item(74, 11).Value = 1444.6111
$ws.Cells.Item(74, 12).Value = 2500
$ws.Cells.Item(74, 13).Value = -570.6111000000001
$ws.Cells.Item(74, 14).Value = -4248
$ws.Cells.Item(77, 8).Value = 1500.1578
$ws.Cells.Item(77, 9).Value = 1444.6111
$ws.Cells.Item(77, 10).Value = 2500
$ws.Cells.Item(77, 11).Value = 7223.0555
$ws.Cells.Item(77, 12).Value = 12500
$ws.Cells.Item(77, 13).Value = -2855.0555
$ws.Cells.Item(77, 14).Value = -21236
$ws.Cells.Item(122, 8).Value = 1222549.5
$ws.Cells.Item(122, 10).Value = 1399.6666
$ws.Cells.Item(122, 12).Value = 4198.9998
$ws.Cells.Item(122, 14).Value = -9098.9998
$ws.Cells.Item(132, 8).Value = 2873.5
$ws.Cells.Item(132, 9).Value = 1523.9615
$ws.Cells.Item(132, 10).Value = 4627.9
$ws.Cells.Item(132, 11).Value = 4571.8845
$ws.Cells.Item(132, 12).Value = 13883.7
$ws.Cells.Item(132, 13).Value = -2041.8845
$ws.Cells.Item(132, 14).Value = -18943.7
$ws.Cells.Item(134, 8).Value = 43009.668
$ws.Cells.Item(134, 10).Value = 43009.668
$ws.Cells.Item(134, 12).Value = 43009.668
$ws.Cells.Item(134, 14).Value = -53149.668
$ws.Cells.Item(136, 8).Value = 9537.643
$ws.Cells.Item(136, 9).Value = 11554.454
$ws.Cells.Item(136, 11).Value = 34663.362
$ws.Cells.Item(136, 13).Value = -32113.362
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3997.3914
$ws.Cells.Item(134, 9).Value = 4331.5137
$ws.Cells.Item(134, 10).Value = 2623.7778
$ws.Cells.Item(134, 11).Value = 12994.5411
$ws.Cells.Item(134, 12).Value = 7871.3334
$ws.Cells.Item(134, 13).Value = -10459.5411
$ws.Cells.Item(134, 14).Value = -12941.3334
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 240.675
$ws.Cells.Item(107, 10).Value = 256.3793
$ws.Cells.Item(107, 12).Value = 256.3793
$ws.Cells.Item(107, 14).Value = -4096.3793
$ws.Cells.Item(122, 8).Value = 1174.2307
$ws.Cells.Item(122, 9).Value = 810
$ws.Cells.Item(122, 10).Value = 1599.1666
$ws.Cells.Item(122, 11).Value = 2430
$ws.Cells.Item(122, 12).Value = 4797.4998
$ws.Cells.Item(122, 13).Value = 20
$ws.Cells.Item(122, 14).Value = -9697.4998
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 6667244
$ws.Cells.Item(113, 10).Value = 1250649.2
$ws.Cells.Item(113, 12).Value = 3751947.6
$ws.Cells.Item(113, 14).Value = -3756287.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 10000
$ws.Cells.Item(10, 9).Value = 5000
$ws.Cells.Item(10, 10).Value = 15000
$ws.Cells.Item(10, 11).Value = 5000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = -4831
$ws.Cells.Item(10, 14).Value = -15338
$ws.Cells.Item(102, 8).Value = 1007.4
$ws.Cells.Item(102, 9).Value = 939.1429000000001
$ws.Cells.Item(102, 10).Value = 1166.6666
$ws.Cells.Item(102, 11).Value = 939.1429000000001
$ws.Cells.Item(102, 12).Value = 1166.6666
$ws.Cells.Item(102, 13).Value = 682.8570999999999
$ws.Cells.Item(102, 14).Value = -4410.6666
$ws.Cells.Item(122, 8).Value = 1966204.8
$ws.Cells.Item(122, 9).Value = 2494582.5
$ws.Cells.Item(122, 10).Value = 3658.2856
$ws.Cells.Item(122, 11).Value = 7483747.5
$ws.Cells.Item(122, 12).Value = 10974.8568
$ws.Cells.Item(122, 13).Value = -7481297.5
$ws.Cells.Item(122, 14).Value = -15874.8568
$ws.Cells.Item(132, 8).Value = 2574.946
$ws.Cells.Item(132, 9).Value = 2350.5293
$ws.Cells.Item(132, 10).Value = 2765.7
$ws.Cells.Item(132, 11).Value = 7051.5879
$ws.Cells.Item(132, 12).Value = 8297.099999999999
$ws.Cells.Item(132, 13).Value = -4521.5879
$ws.Cells.Item(132, 14).Value = -13357.1
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 71430500
$ws.Cells.Item(40, 9).Value = 111112710
$ws.Cells.Item(40, 10).Value = 2522
$ws.Cells.Item(40, 11).Value = 111112710
$ws.Cells.Item(40, 12).Value = 2522
$ws.Cells.Item(40, 13).Value = -111112574
$ws.Cells.Item(40, 14).Value = -2794
$ws.Cells.Item(122, 8).Value = 3133506
$ws.Cells.Item(122, 9).Value = 4465798.5
$ws.Cells.Item(122, 11).Value = 13397395.5
$ws.Cells.Item(122, 13).Value = -13394945.5
$ws.Cells.Item(132, 8).Value = 8188477.5
$ws.Cells.Item(132, 9).Value = 13893953
$ws.Cells.Item(132, 10).Value = 2360.9565
$ws.Cells.Item(132, 11).Value = 41681859
$ws.Cells.Item(132, 12).Value = 7082.869499999999
$ws.Cells.Item(132, 13).Value = -41679329
$ws.Cells.Item(132, 14).Value = -12142.8695
$ws.Cells.Item(134, 8).Value = 44750
$ws.Cells.Item(134, 10).Value = 44750
$ws.Cells.Item(134, 12).Value = 44750
$ws.Cells.Item(134, 14).Value = -54890
$ws.Cells.Item(136, 8).Value = 9321.517
$ws.Cells.Item(136, 9).Value = 7511.727
$ws.Cells.Item(136, 10).Value = 13745.444
$ws.Cells.Item(136, 11).Value = 22535.181
$ws.Cells.Item(136, 12).Value = 41236.33199999999
$ws.Cells.Item(136, 13).Value = -19985.181
$ws.Cells.Item(136, 14).Value = -46336.33199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1334.9166
$ws.Cells.Item(122, 9).Value = 1252.375
$ws.Cells.Item(122, 11).Value = 3757.125
$ws.Cells.Item(122, 13).Value = -1307.125
$ws.Cells.Item(132, 8).Value = 1058.7727
$ws.Cells.Item(132, 9).Value = 831.4146
$ws.Cells.Item(132, 10).Value = 4166
$ws.Cells.Item(132, 11).Value = 2494.2438
$ws.Cells.Item(132, 12).Value = 12498
$ws.Cells.Item(132, 13).Value = 35.75620000000026
$ws.Cells.Item(132, 14).Value = -17558
$ws.Cells.Item(136, 8).Value = 1319.6207
$ws.Cells.Item(136, 9).Value = 686.5625
$ws.Cells.Item(136, 10).Value = 2098.7693
$ws.Cells.Item(136, 11).Value = 2059.6875
$ws.Cells.Item(136, 12).Value = 6296.3079
$ws.Cells.Item(136, 13).Value = 490.3125
$ws.Cells.Item(136, 14).Value = -11396.3079